$wb = $excel.ActiveWorkbook

# Sheet "OFF" (sheet1) - row 3 (A3 = "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 426
$wsOff.Range("C3").Value = 272
$wsOff.Range("D3").Value = 74
$wsOff.Range("E3").Value = 35
$wsOff.Range("F3").Value = 9

# Sheet "DEF" (sheet2) - row 3 (A3 = "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 438
$wsDef.Range("C3").Value = 333
$wsDef.Range("D3").Value = 83
$wsDef.Range("E3").Value = 36
